# Update "想去人数" (want-to-go count) figures in 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3746
$ws1.Range("F5").Value  = 2267
$ws1.Range("F6").Value  = 444
$ws1.Range("F11").Value = 84
$ws1.Range("F12").Value = 1389
$ws1.Range("F14").Value = 2233
$ws1.Range("F15").Value = 160

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 3746
$ws4.Range("F5").Value  = 2267
$ws4.Range("F6").Value  = 444
$ws4.Range("F12").Value = 84
$ws4.Range("F15").Value = 1389
$ws4.Range("F17").Value = 2234
$ws4.Range("F18").Value = 160
